# Auto-generated Excel COM-interop script
# Applies "csr search remittance cases" edits to the SQL worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")
$ws.Activate()

$ws.Range("B42").Value = "select PROV_TAX_ID_NBR`nfrom PP001.PROVIDER as p,PP001.CONSOLIDATED_PAYMENT cp,OLE.PROC_CTL pc`nwhere p.PROV_KEY_ID =cp.PROV_KEY_ID`nand cp.PROC_CTL_ID=pc.PROC_CTL_ID`nand cp.setl_dt between (current date - 60 days) and current date`nand pc.EXTRACT_STS_CD='C' `nfetch first row only with ur"
$ws.Range("C42").Value = "Get prov tin number (DOP)"

$ws.Range("B43").Value = "Select c.PTNT_ACCT_NBR, cp.SETL_DT`nfrom PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup,OLE.PROC_CTL pc`nwhere cp.prov_key_id = p.prov_key_id`nand cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID`nand cup.CLM_KEY_ID = c.CLM_KEY_ID and cp.PROC_CTL_ID=pc.PROC_CTL_ID`nand p.PROV_TAX_ID_NBR ='{`$tin}'`nand pc.EXTRACT_STS_CD='C'`nand cp.SETL_DT between (current date - 60 days) and current date `norder by cp.SETL_DT DESC`nfetch first row only"

$ws.Range("B44").Value = "Select sr.SBSCR_ID, cp.SETL_DT`nfrom PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup, PP001.SUBSCRIBER sr,OLE.PROC_CTL pc`nwhere cp.prov_key_id = p.prov_key_id`nand cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID`nand cup.CLM_KEY_ID = c.CLM_KEY_ID`nand c.SBSCR_KEY_ID = sr.SBSCR_KEY_ID and cp.PROC_CTL_ID=pc.PROC_CTL_ID and pc.EXTRACT_STS_CD='C'`nand p.PROV_TAX_ID_NBR ='{`$tin}'`nand cp.SETL_DT between (current date - 60 days) and current date  `norder by cp.SETL_DT DESC`nfetch first row only"

$ws.Range("B45").Value = "Select p.PROV_NPI_NBR, cp.SETL_DT`nfrom PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p,OLE.PROC_CTL pc`nwhere cp.prov_key_id = p.prov_key_id and cp.PROC_CTL_ID=pc.PROC_CTL_ID`nand p.PROV_TAX_ID_NBR ='{`$tin}'  and pc.EXTRACT_STS_CD='C'`nand cp.SETL_DT between (current date - 60 days) and current date `norder by cp.SETL_DT DESC`nfetch first row only with ur"

$ws.Range("B46").Value = "Select c.CLM_NBR, cp.SETL_DT`nfrom PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup,OLE.PROC_CTL pc`nwhere cp.prov_key_id = p.prov_key_id`nand cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID`nand cup.CLM_KEY_ID = c.CLM_KEY_ID and cp.PROC_CTL_ID=pc.PROC_CTL_ID and pc.EXTRACT_STS_CD='C'`nand p.PROV_TAX_ID_NBR = '{`$tin}'`nand cp.SETL_DT between (current date - 60 days) and current date `norder by cp.SETL_DT DESC`nfetch first row only with ur"

$ws.Range("B47").Value = "Select c.PTNT_FST_NM, c.PTNT_LST_NM, cp.SETL_DT`nfrom PP001.CONSOLIDATED_PAYMENT cp, PP001.PROVIDER p, PP001.CLAIM c, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CLAIM_UNCONSOLIDATED_PAYMENT cup,OLE.PROC_CTL pc`nwhere cp.prov_key_id = p.prov_key_id`nand cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID`nand cup.CLM_KEY_ID = c.CLM_KEY_ID and cp.PROC_CTL_ID=pc.PROC_CTL_ID and pc.EXTRACT_STS_CD='C'`nand p.PROV_TAX_ID_NBR = '{`$tin}'`nand cp.SETL_DT between (current date - 60 days) and current date  and c.PTNT_FST_NM <> ''`nand c.PTNT_LST_NM <> ''`norder by cp.SETL_DT DESC`nfetch first row only with ur"
$ws.Range("C47").Value = "Get Patient first name last name for a tin"

$ws.Range("B48").Value = "Select p.PROV_TAX_ID_NBR,cp.SETL_DT`nfrom PP001.PROVIDER p, PP001.UNCONSOLIDATED_PAYMENT ucp, PP001.CONSOLIDATED_PAYMENT cp,`nPP001.CLAIM_UNCONSOLIDATED_PAYMENT cup,OLE.PROC_CTL pc`nwhere ucp.prov_key_id = p.prov_key_id`nand ucp.UCONSL_PAY_KEY_ID = cup.UCONSL_PAY_KEY_ID `nand cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand cp.PROC_CTL_ID=pc.PROC_CTL_ID and pc.EXTRACT_STS_CD='C'`nand cup.CLM_PAY_AMT = '0.00'`nand ucp.PROC_DTTM  between (current date - 180 days) and current date `norder by ucp.PROC_DTTM DESC`nfetch first row only with ur"
$ws.Range("C48").Value = "To get tin by DOP and Zero Payment claims- query taking more than 2 minutes ..need to optimized"

$ws.Range("B49").Value = "select p.PROV_TAX_ID_NBR `nfrom PP001.PROVIDER p, PP001.UNCONSOLIDATED_PAYMENT ucp, OLE.PROC_CTL pc, PP001.CONSOLIDATED_PAYMENT cp`nwhere cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR and p.PROV_KEY_ID =ucp.PROV_KEY_ID `nand cp.PROC_CTL_ID=pc.PROC_CTL_ID `nand pc.EXTRACT_STS_CD='C'`nand cp.setl_dt between current date - 6 MONTHS and current date`nfetch first row only with ur"
$ws.Range("C49").Value = "Get the tin No for Electronic Payment No Search "

$ws.Range("B50").Value = "select cp.DSPL_CONSL_PAY_NBR`nfrom PP001.PROVIDER p, PP001.UNCONSOLIDATED_PAYMENT ucp,PP001.CONSOLIDATED_PAYMENT cp,OLE.PROC_CTL pc`nwhere cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand p.PROV_KEY_ID =ucp.PROV_KEY_ID `nand cp.PROC_CTL_ID=pc.PROC_CTL_ID `nand p.PROV_TAX_ID_NBR='{`$tin}'`nand pc.EXTRACT_STS_CD='C'`nand cp.setl_dt between current date - 6 MONTHS and current date`nfetch first row only with ur"
$ws.Range("C50").Value = "Get the Electronic No For a Tin from above query"

$ws.Range("B51").Value = "select p.PROV_TAX_ID_NBR `nfrom PP001.PROVIDER p, PP001.UNCONSOLIDATED_PAYMENT ucp, OLE.PROC_CTL pc, PP001.CONSOLIDATED_PAYMENT cp`nwhere cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR and p.PROV_KEY_ID =ucp.PROV_KEY_ID `nand cp.PROC_CTL_ID=pc.PROC_CTL_ID `nand pc.EXTRACT_STS_CD='C'`nand cp.setl_dt between current date - 6 MONTHS and current date`nfetch first row only with ur"
$ws.Range("C51").Value = "Get the tin No for Check No Search & Electronic No"

$ws.Range("B52").Value = "select ucp.UCONSL_PAY_NBR`nfrom PP001.PROVIDER p, PP001.UNCONSOLIDATED_PAYMENT ucp,PP001.CONSOLIDATED_PAYMENT cp,OLE.PROC_CTL pc`nwhere cp.CONSL_PAY_NBR = ucp.CONSL_PAY_NBR`nand p.PROV_KEY_ID =ucp.PROV_KEY_ID `nand cp.PROC_CTL_ID=pc.PROC_CTL_ID `nand p.PROV_TAX_ID_NBR='{`$tin}'`nand pc.EXTRACT_STS_CD='C'`nand cp.setl_dt between current date - 6 MONTHS and current date`nfetch first row only with ur"
$ws.Range("C52").Value = "Get the Check No for a tin from above query"

# Update the view: scroll to show rows starting at A47, select B48
$ws.Application.ActiveWindow.ScrollRow = 47
$ws.Range("B48").Select()
